# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 545
$ws1.Range("F7").Value = 2732
$ws1.Range("F8").Value = 462
$ws1.Range("F9").Value = 7634
$ws1.Range("F11").Value = 465
$ws1.Range("F13").Value = 292

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 545
$ws4.Range("F9").Value = 2732
$ws4.Range("F10").Value = 462
$ws4.Range("F11").Value = 7634
$ws4.Range("F13").Value = 465
$ws4.Range("F17").Value = 292
